$wb = $excel.ActiveWorkbook

# =========================================================================
# 1) Insert the new "2022-Q4" worksheet right before "2022-Q3"
# =========================================================================
$refSheet = $wb.Worksheets.Item("2022-Q3")
$newSheetRaw = $wb.Worksheets.Add($refSheet)
$newSheetRaw.Name = "2022-Q4"

# Re-fetch sheet handles by name -- the object returned directly by Add()
# does not reliably carry cross-sheet style (PasteSpecial) writes in this
# runtime, so all further work happens through freshly-looked-up handles.
$src = $wb.Worksheets.Item("2022-Q3")
$ws = $wb.Worksheets.Item("2022-Q4")

# ---- Data for new "2022-Q4" sheet ----
$data = @(
  @("481010","工银中小盘混合","15.71","91.30","3.88","0.6095","3"),
  @("050004","博时精选混合A","18.58","72.42","2.96","0.5500","7"),
  @("012985","平安优势回报1年持有混合A","12.27","94.40","3.97","0.4871","5"),
  @("001236","博时丝路主题股票A","10.52","87.16","2.73","0.2872","9"),
  @("012917","平安优势领航1年持有期混合A","6.72","93.97","3.71","0.2493","7"),
  @("013417","博时核心资产精选混合A","6.96","80.13","3.37","0.2346","6"),
  @("210003","金鹰行业优势混合","5.27","86.97","4.13","0.2177","8"),
  @("002450","平安睿享文娱灵活配置混合A","3.85","94.03","4.24","0.1632","5"),
  @("010126","平安价值成长混合A","3.21","94.03","4.06","0.1303","5"),
  @("013365","汇添富产业升级混合A","4.03","85.05","3.21","0.1294","10"),
  @("014062","景顺长城专精特新量化优选股票A","7.66","91.15","1.51","0.1157","9"),
  @("011828","平安睿享成长混合A","2.42","92.96","4.52","0.1094","5"),
  @("006101","平安优势产业灵活配置混合C","2.85","94.90","3.50","0.0998","9"),
  @("002451","平安睿享文娱灵活配置混合C","1.97","94.03","4.24","0.0835","5"),
  @("006100","平安优势产业灵活配置混合A","2.15","94.90","3.50","0.0752","9"),
  @("013711","广发成长新动能混合C","2.39","84.33","3.10","0.0741","10"),
  @("162717","广发成长新动能混合A","2.33","84.33","3.10","0.0722","10"),
  @("013687","平安成长龙头1年持有混合A","1.21","94.87","5.03","0.0609","4"),
  @("014063","景顺长城专精特新量化优选股票C","3.88","91.15","1.51","0.0586","9"),
  @("010127","平安价值成长混合C","1.33","94.03","4.06","0.0540","5"),
  @("005265","博时厚泽回报灵活配置混合A","1.63","77.88","3.04","0.0496","7"),
  @("011829","平安睿享成长混合C","1.09","92.96","4.52","0.0493","5"),
  @("012986","平安优势回报1年持有混合C","1.14","94.40","3.97","0.0453","5"),
  @("005266","博时厚泽回报灵活配置混合C","0.90","77.88","3.04","0.0274","7"),
  @("013688","平安成长龙头1年持有混合C","0.51","94.87","5.03","0.0257","4"),
  @("002556","博时丝路主题股票C","0.92","87.16","2.73","0.0251","9"),
  @("011340","博时战略新材料主题混合A","0.84","79.48","2.95","0.0248","6"),
  @("007894","平安估值精选混合C","0.50","94.39","4.32","0.0216","5"),
  @("011341","博时战略新材料主题混合C","0.56","79.48","2.95","0.0165","6"),
  @("014212","博时研究优享混合A","0.59","79.20","2.69","0.0159","10"),
  @("007893","平安估值精选混合A","0.34","94.39","4.32","0.0147","5"),
  @("013418","博时核心资产精选混合C","0.38","80.13","3.37","0.0128","6"),
  @("016370","信澳业绩驱动混合A","0.77","30.31","1.66","0.0128","10"),
  @("013366","汇添富产业升级混合C","0.28","85.05","3.21","0.0090","10"),
  @("016371","信澳业绩驱动混合C","0.25","30.31","1.66","0.0042","10"),
  @("012918","平安优势领航1年持有期混合C","0.08","93.97","3.71","0.0030","7"),
  @("014213","博时研究优享混合C","0.08","79.20","2.69","0.0022","10"),
  @("016751","博时精选混合C","0.00","72.42","2.96","0","7")
)

# ---- Header row (row 1), columns B..H ----
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# ---- Column A (row index, 0-based), numeric ----
for ($i = 0; $i -lt $data.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i
}

# ---- Columns B..G, forced to text so numeric-looking strings (fund codes
#      with leading zeros, percentages, NAV figures, ...) stay text, which
#      is how the source data is encoded ----
$lastRow = $data.Length + 1
$textRange = $ws.Range("B2:G$lastRow")
$textRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $row[0]   # 基金代码
    $ws.Cells.Item($r, 3).Value = $row[1]   # 基金名称
    $ws.Cells.Item($r, 4).Value = $row[2]   # 基金规模
    $ws.Cells.Item($r, 5).Value = $row[3]   # 股票总仓位
    $ws.Cells.Item($r, 6).Value = $row[4]   # 仓位占比
    $ws.Cells.Item($r, 7).Value = $row[5]   # 持有市值(亿元)
    $ws.Cells.Item($r, 8).Value = [double]$row[6]   # 仓位排名 (number)
}

# Remove the temporary "@" text format so the cells end up with the default
# style (matching the source file, which carries no special style on the
# data cells) while keeping the values typed as text.
$textRange.ClearFormats()

# The very last row's "持有市值(亿元)" is stored as a genuine number (0) in
# the source data rather than text -- fix it up after the ClearFormats pass.
$ws.Cells.Item($lastRow, 7).Value = 0

# ---- Styles: reuse the existing bold/border/center style (s=2) from the
#      "2022-Q3" sheet's header row and index column, via copy/paste of
#      formats only (-4122 = xlPasteFormats) ----
$src.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

$src.Range("A2").Copy()
$ws.Range("A2:A$lastRow").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# =========================================================================
# 2) Update the "总计" summary sheet: insert a new row 2 for "2022-Q4" and
#    shift the existing quarters down, renumbering the index column.
# =========================================================================
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 38
$total.Cells.Item(2, 4).Value = 4.22

for ($r = 3; $r -le 7; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

# Row 2 picked up row 1's (header) formatting on the Insert; fix it to match
# the other data rows: B2:D2 unstyled, A2 reusing the A-column's s=2 style.
$total.Range("B2:D2").ClearFormats()
$total.Cells.Item(3, 1).Copy()
$total.Cells.Item(2, 1).PasteSpecial(-4122)
$excel.CutCopyMode = 0

Write-Host "edit complete"
